$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model results")

# --- Row 2: Clavigralla shadabi Benin ---
# Remove r.model.h (G2) and r.model.f (H2); add a "No diurnal variation" note
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = "No diurnal variation"

# --- Row 3: Clavigralla tomentosicollis Benin ---
# Re-fit r.TPC.h / r.TPC.f; remove r.model.h / r.model.f; add note
$ws.Range("E3").Value = 0.13600000000000001
$ws.Range("F3").Value = 0.153
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = "No diurnal variation"

# --- Row 4: Clavigralla tomentosicollis Burkina Faso ---
# Re-fit r.TPC.h / r.TPC.f; remove r.model.h; clear r.model.f value (keep its format); add note
$ws.Range("E4").Value = 0.15
$ws.Range("E4").NumberFormat = "0.000"
$ws.Range("F4").Value = 0.153
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").Value = "No diurnal variation"

# --- Row 5: Apolygus lucorum China Dafeng ---
# Fill in newly-fit parameters and note
$ws.Range("E5").Value = 0.20599999999999999
$ws.Range("F5").Value = 0.214
$ws.Range("G5").Value = 0.048
$ws.Range("G5").NumberFormat = "0.0000"
$ws.Range("H5").Value = 0.042
$ws.Range("H5").NumberFormat = "0.0000"
$ws.Range("I5").Value = "No diurnal variation"

# --- Column widths: best-fit widths for the new "subfamily" (D) and "Notes" (I) columns ---
$ws.Columns.Item(4).ColumnWidth = 12.6666666666667
$ws.Columns.Item(9).ColumnWidth = 16.6666666666667

# Minor re-fit of column A width
$ws.Columns.Item(1).ColumnWidth = 34.5

# --- Cursor/selection moved to I11 ---
$ws.Range("I11").Select()
